$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.056197333333333
$ws.Range("H2").Value = 3.168592
$ws.Range("I2").Value = 0.01247237710445079
$ws.Range("J2").Value = 0.01398563433468744
$ws.Range("M2").Value = 14.23612
$ws.Range("N2").Value = 42.70836
$ws.Range("O2").Value = 0.07600219901258977
$ws.Range("P2").Value = 0.09845490883293605
$ws.Range("Q2").Value = 15.03615198101333
$ws.Range("R2").Value = 135.32536782912
$ws.Range("S2").Value = 0.0009479280868525372
$ws.Range("T2").Value = 0.001376954353392433

$ws.Range("G3").Value = 1.056197333333333
$ws.Range("H3").Value = 3.168592
$ws.Range("I3").Value = 0.01247237710445079
$ws.Range("J3").Value = 0.01398563433468744
$ws.Range("O3").Value = 0.1644984167819043
$ws.Range("P3").Value = 0.2130948424892534
$ws.Range("Q3").Value = 32.54410040108266
$ws.Range("R3").Value = 292.896903609744
$ws.Range("S3").Value = 0.002051686287189027
$ws.Range("T3").Value = 0.002980266545662515

$ws.Range("G4").Value = 1.056197333333333
$ws.Range("H4").Value = 3.168592
$ws.Range("I4").Value = 0.01247237710445079
$ws.Range("J4").Value = 0.01398563433468744
$ws.Range("M4").Value = 5.407681
$ws.Range("N4").Value = 16.223043
$ws.Range("O4").Value = 0.02886992014387351
$ws.Range("P4").Value = 0.03739872520410058
$ws.Range("Q4").Value = 5.711578251717333
$ws.Range("R4").Value = 51.404204265456
$ws.Range("S4").Value = 0.0003600765310097706
$ws.Range("T4").Value = 0.0005230448952880098

$ws.Range("G5").Value = 1.056197333333333
$ws.Range("H5").Value = 3.168592
$ws.Range("I5").Value = 0.01247237710445079
$ws.Range("J5").Value = 0.01398563433468744
$ws.Range("M5").Value = 128.149857
$ws.Range("N5").Value = 256.299714
$ws.Range("O5").Value = 0.6841520677789258
$ws.Range("P5").Value = 0.5908436890523912
$ws.Range("Q5").Value = 135.351537230448
$ws.Range("R5").Value = 812.1092233826879
$ws.Range("S5").Value = 0.008533002586128539
$ws.Range("T5").Value = 0.008263323784044515

$ws.Range("G6").Value = 1.056197333333333
$ws.Range("H6").Value = 3.168592
$ws.Range("I6").Value = 0.01247237710445079
$ws.Range("J6").Value = 0.01398563433468744
$ws.Range("M6").Value = 8.705771666666665
$ws.Range("N6").Value = 26.117315
$ws.Range("O6").Value = 0.04647739628270661
$ws.Range("P6").Value = 0.06020783442131875
$ws.Range("Q6").Value = 9.195012818942219
$ws.Range("R6").Value = 82.75511537047998
$ws.Range("S6").Value = 0.0005796836132709163
$ws.Range("T6").Value = 0.000842044756299972

$ws.Range("I7").Value = 0.6620593097549599
$ws.Range("J7").Value = 0.7423861014276285
$ws.Range("M7").Value = 14.23612
$ws.Range("N7").Value = 42.70836
$ws.Range("O7").Value = 0.07600219901258977
$ws.Range("P7").Value = 0.09845490883293605
$ws.Range("Q7").Value = 798.1497286806667
$ws.Range("R7").Value = 7183.347558126001
$ws.Range("S7").Value = 0.05031796341813428
$ws.Range("T7").Value = 0.07309155593489598

$ws.Range("I8").Value = 0.6620593097549599
$ws.Range("J8").Value = 0.7423861014276285
$ws.Range("O8").Value = 0.1644984167819043
$ws.Range("P8").Value = 0.2130948424892534
$ws.Range("S8").Value = 0.1089077082704113
$ws.Range("T8").Value = 0.1581986493499314

$ws.Range("I9").Value = 0.6620593097549599
$ws.Range("J9").Value = 0.7423861014276285
$ws.Range("M9").Value = 5.407681
$ws.Range("N9").Value = 16.223043
$ws.Range("O9").Value = 0.02886992014387351
$ws.Range("P9").Value = 0.03739872520410058
$ws.Range("Q9").Value = 303.1822661611167
$ws.Range("R9").Value = 2728.640395450051
$ws.Range("S9").Value = 0.01911359940313371
$ws.Range("T9").Value = 0.02776429380263542

$ws.Range("I10").Value = 0.6620593097549599
$ws.Range("J10").Value = 0.7423861014276285
$ws.Range("M10").Value = 128.149857
$ws.Range("N10").Value = 256.299714
$ws.Range("O10").Value = 0.6841520677789258
$ws.Range("P10").Value = 0.5908436890523912
$ws.Range("Q10").Value = 7184.736683521651
$ws.Range("R10").Value = 43108.4201011299
$ws.Range("S10").Value = 0.4529492457611441
$ws.Range("T10").Value = 0.4386341428687227

$ws.Range("I11").Value = 0.6620593097549599
$ws.Range("J11").Value = 0.7423861014276285
$ws.Range("M11").Value = 8.705771666666665
$ws.Range("N11").Value = 26.117315
$ws.Range("O11").Value = 0.04647739628270661
$ws.Range("P11").Value = 0.06020783442131875
$ws.Range("Q11").Value = 488.0901041650278
$ws.Range("R11").Value = 4392.81093748525
$ws.Range("S11").Value = 0.03077079290213648
$ws.Range("T11").Value = 0.044697459471443

$ws.Range("G12").Value = 0.04559766666666667
$ws.Range("H12").Value = 0.136793
$ws.Range("I12").Value = 0.0005384517417354892
$ws.Range("J12").Value = 0.0006037813885615125
$ws.Range("M12").Value = 14.23612
$ws.Range("N12").Value = 42.70836
$ws.Range("O12").Value = 0.07600219901258977
$ws.Range("P12").Value = 0.09845490883293605
$ws.Range("Q12").Value = 0.6491338543866667
$ws.Range("R12").Value = 5.84220468948
$ws.Range("S12").Value = 0.00004092351643405624
$ws.Range("T12").Value = 0.00005944524156584725

$ws.Range("G13").Value = 0.04559766666666667
$ws.Range("H13").Value = 0.136793
$ws.Range("I13").Value = 0.0005384517417354892
$ws.Range("J13").Value = 0.0006037813885615125
$ws.Range("O13").Value = 0.1644984167819043
$ws.Range("P13").Value = 0.2130948424892534
$ws.Range("Q13").Value = 1.404978970522333
$ws.Range("R13").Value = 12.644810734701
$ws.Range("S13").Value = 0.00008857445902894681
$ws.Range("T13").Value = 0.0001286626998934582

$ws.Range("G14").Value = 0.04559766666666667
$ws.Range("H14").Value = 0.136793
$ws.Range("I14").Value = 0.0005384517417354892
$ws.Range("J14").Value = 0.0006037813885615125
$ws.Range("M14").Value = 5.407681
$ws.Range("N14").Value = 16.223043
$ws.Range("O14").Value = 0.02886992014387351
$ws.Range("P14").Value = 0.03739872520410058
$ws.Range("Q14").Value = 0.2465776356776667
$ws.Range("R14").Value = 2.219198721099
$ws.Range("S14").Value = 0.00001554505878523318
$ws.Range("T14").Value = 0.00002258065423416229

$ws.Range("G15").Value = 0.04559766666666667
$ws.Range("H15").Value = 0.136793
$ws.Range("I15").Value = 0.0005384517417354892
$ws.Range("J15").Value = 0.0006037813885615125
$ws.Range("M15").Value = 128.149857
$ws.Range("N15").Value = 256.299714
$ws.Range("O15").Value = 0.6841520677789258
$ws.Range("P15").Value = 0.5908436890523912
$ws.Range("Q15").Value = 5.843334462867
$ws.Range("R15").Value = 35.060006777202
$ws.Range("S15").Value = 0.0003683828725074991
$ws.Range("T15").Value = 0.0003567404229988592

$ws.Range("G16").Value = 0.04559766666666667
$ws.Range("H16").Value = 0.136793
$ws.Range("I16").Value = 0.0005384517417354892
$ws.Range("J16").Value = 0.0006037813885615125
$ws.Range("M16").Value = 8.705771666666665
$ws.Range("N16").Value = 26.117315
$ws.Range("O16").Value = 0.04647739628270661
$ws.Range("P16").Value = 0.06020783442131875
$ws.Range("Q16").Value = 0.3969628745327777
$ws.Range("R16").Value = 3.572665870795
$ws.Range("S16").Value = 0.00002502583497975393
$ws.Range("T16").Value = 0.00003635236986918546

$ws.Range("G17").Value = 27.488287
$ws.Range("H17").Value = 54.976574
$ws.Range("I17").Value = 0.3246024872429512
$ws.Range("J17").Value = 0.2426573888143015
$ws.Range("M17").Value = 14.23612
$ws.Range("N17").Value = 42.70836
$ws.Range("O17").Value = 0.07600219901258977
$ws.Range("P17").Value = 0.09845490883293605
$ws.Range("Q17").Value = 391.32655232644
$ws.Range("R17").Value = 2347.95931395864
$ws.Range("S17").Value = 0.02467050283542041
$ws.Range("T17").Value = 0.02389081109335037

$ws.Range("G18").Value = 27.488287
$ws.Range("H18").Value = 54.976574
$ws.Range("I18").Value = 0.3246024872429512
$ws.Range("J18").Value = 0.2426573888143015
$ws.Range("O18").Value = 0.1644984167819043
$ws.Range("P18").Value = 0.2130948424892534
$ws.Range("Q18").Value = 846.983365464953
$ws.Range("R18").Value = 5081.900192789718
$ws.Range("S18").Value = 0.05339659523493375
$ws.Range("T18").Value = 0.05170903804823708

$ws.Range("G19").Value = 27.488287
$ws.Range("H19").Value = 54.976574
$ws.Range("I19").Value = 0.3246024872429512
$ws.Range("J19").Value = 0.2426573888143015
$ws.Range("M19").Value = 5.407681
$ws.Range("N19").Value = 16.223043
$ws.Range("O19").Value = 0.02886992014387351
$ws.Range("P19").Value = 0.03739872520410058
$ws.Range("Q19").Value = 148.647887332447
$ws.Range("R19").Value = 891.887323994682
$ws.Range("S19").Value = 0.009371247885206718
$ws.Range("T19").Value = 0.009075077003010652

$ws.Range("G20").Value = 27.488287
$ws.Range("H20").Value = 54.976574
$ws.Range("I20").Value = 0.3246024872429512
$ws.Range("J20").Value = 0.2426573888143015
$ws.Range("M20").Value = 128.149857
$ws.Range("N20").Value = 256.299714
$ws.Range("O20").Value = 0.6841520677789258
$ws.Range("P20").Value = 0.5908436890523912
$ws.Range("Q20").Value = 3522.620048224959
$ws.Range("R20").Value = 14090.48019289983
$ws.Range("S20").Value = 0.2220774628534474
$ws.Range("T20").Value = 0.1433725867828623

$ws.Range("G21").Value = 27.488287
$ws.Range("H21").Value = 54.976574
$ws.Range("I21").Value = 0.3246024872429512
$ws.Range("J21").Value = 0.2426573888143015
$ws.Range("M21").Value = 8.705771666666665
$ws.Range("N21").Value = 26.117315
$ws.Range("O21").Value = 0.04647739628270661
$ws.Range("P21").Value = 0.06020783442131875
$ws.Range("Q21").Value = 239.3067501298016
$ws.Range("R21").Value = 1435.84050077881
$ws.Range("S21").Value = 0.01508667843394286
$ws.Range("T21").Value = 0.01460987588684103

$ws.Range("G22").Value = 0.027723
$ws.Range("H22").Value = 0.08316900000000001
$ws.Range("I22").Value = 0.0003273741559027063
$ws.Range("J22").Value = 0.0003670940348210247
$ws.Range("M22").Value = 14.23612
$ws.Range("N22").Value = 42.70836
$ws.Range("O22").Value = 0.07600219901258977
$ws.Range("P22").Value = 0.09845490883293605
$ws.Range("Q22").Value = 0.39466795476
$ws.Range("R22").Value = 3.55201159284
$ws.Range("S22").Value = 0.00002488115574849607
$ws.Range("T22").Value = 0.00003614220973141863

$ws.Range("G23").Value = 0.027723
$ws.Range("H23").Value = 0.08316900000000001
$ws.Range("I23").Value = 0.0003273741559027063
$ws.Range("J23").Value = 0.0003670940348210247
$ws.Range("O23").Value = 0.1644984167819043
$ws.Range("P23").Value = 0.2130948424892534
$ws.Range("Q23").Value = 0.854215464237
$ws.Range("R23").Value = 7.687939178133001
$ws.Range("S23").Value = 0.0000538525303413075
$ws.Range("T23").Value = 0.00007822584552893074

$ws.Range("G24").Value = 0.027723
$ws.Range("H24").Value = 0.08316900000000001
$ws.Range("I24").Value = 0.0003273741559027063
$ws.Range("J24").Value = 0.0003670940348210247
$ws.Range("M24").Value = 5.407681
$ws.Range("N24").Value = 16.223043
$ws.Range("O24").Value = 0.02886992014387351
$ws.Range("P24").Value = 0.03739872520410058
$ws.Range("Q24").Value = 0.149917140363
$ws.Range("R24").Value = 1.349254263267
$ws.Range("S24").Value = 0.000009451265738079126
$ws.Range("T24").Value = 0.00001372884893233603

$ws.Range("G25").Value = 0.027723
$ws.Range("H25").Value = 0.08316900000000001
$ws.Range("I25").Value = 0.0003273741559027063
$ws.Range("J25").Value = 0.0003670940348210247
$ws.Range("M25").Value = 128.149857
$ws.Range("N25").Value = 256.299714
$ws.Range("O25").Value = 0.6841520677789258
$ws.Range("P25").Value = 0.5908436890523912
$ws.Range("Q25").Value = 3.552698485611
$ws.Range("R25").Value = 21.316190913666
$ws.Range("S25").Value = 0.0002239737056982169
$ws.Range("T25").Value = 0.0002168951937627812

$ws.Range("G26").Value = 0.027723
$ws.Range("H26").Value = 0.08316900000000001
$ws.Range("I26").Value = 0.0003273741559027063
$ws.Range("J26").Value = 0.0003670940348210247
$ws.Range("M26").Value = 8.705771666666665
$ws.Range("N26").Value = 26.117315
$ws.Range("O26").Value = 0.04647739628270661
$ws.Range("P26").Value = 0.06020783442131875
$ws.Range("Q26").Value = 0.241350107915
$ws.Range("R26").Value = 2.172150971235
$ws.Range("S26").Value = 0.00001521549837660666
$ws.Range("T26").Value = 0.00002210193686555807
